$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$xlPasteFormats = -4122
$xlRight = -4152

# ====================================================================
# STEP 0: cache existing numeric/formula content of rows 4-6 (these
# rows get rotated: old row6 -> new row4, old row4 -> new row5,
# old row5 -> new row6), and stage copies of every distinct cell style
# that will be needed later, into scratch cells in column N (far away
# from any used column), BEFORE any destination cell is overwritten.
# ====================================================================
$row4vals = @{}
$row5vals = @{}
$row6vals = @{}
foreach ($col in @("D","E","F","G","H","I","J")) {
    $row4vals[$col] = $ws.Range($col + "4").Value2
    $row5vals[$col] = $ws.Range($col + "5").Value2
    $row6vals[$col] = $ws.Range($col + "6").Value2
}

# style staging: style-id -> scratch cell holding that format
$stage = @{
    "9"  = "N1"   # header year style (D3)
    "13" = "N2"   # top data-row label style (A4)
    "14" = "N3"   # middle data-row label style (A5)
    "15" = "N4"   # bottom (bordered) data-row label style (A6)
    "16" = "N5"   # bottom bordered numeric (no decimals shown) style (D6)
    "17" = "N6"   # bottom bordered numeric style (E6)
    "18" = "N7"   # J4 numeric style
    "19" = "N8"   # J5 numeric style
    "20" = "N9"   # J6 bottom bordered numeric style
    "21" = "N10"  # D4 numeric style
    "22" = "N11"  # generic numeric (no border) style
}

$ws.Range("D3").Copy();  $ws.Range($stage["9"]).PasteSpecial($xlPasteFormats)  | Out-Null
$ws.Range("A4").Copy();  $ws.Range($stage["13"]).PasteSpecial($xlPasteFormats) | Out-Null
$ws.Range("A5").Copy();  $ws.Range($stage["14"]).PasteSpecial($xlPasteFormats) | Out-Null
$ws.Range("A6").Copy();  $ws.Range($stage["15"]).PasteSpecial($xlPasteFormats) | Out-Null
$ws.Range("D6").Copy();  $ws.Range($stage["16"]).PasteSpecial($xlPasteFormats) | Out-Null
$ws.Range("E6").Copy();  $ws.Range($stage["17"]).PasteSpecial($xlPasteFormats) | Out-Null
$ws.Range("J4").Copy();  $ws.Range($stage["18"]).PasteSpecial($xlPasteFormats) | Out-Null
$ws.Range("J5").Copy();  $ws.Range($stage["19"]).PasteSpecial($xlPasteFormats) | Out-Null
$ws.Range("J6").Copy();  $ws.Range($stage["20"]).PasteSpecial($xlPasteFormats) | Out-Null
$ws.Range("D4").Copy();  $ws.Range($stage["21"]).PasteSpecial($xlPasteFormats) | Out-Null
$ws.Range("E4").Copy();  $ws.Range($stage["22"]).PasteSpecial($xlPasteFormats) | Out-Null

# style 23 = style 22 (numFmt 164 / vertical center / no border) plus
# horizontal=right alignment -- this is the one genuinely NEW style.
$ws.Range($stage["22"]).Copy(); $ws.Range("N12").PasteSpecial($xlPasteFormats) | Out-Null
$ws.Range("N12").HorizontalAlignment = $xlRight
$stage["23"] = "N12"

function Apply-Style($rangeAddr, $styleId) {
    $ws.Range($stage[$styleId]).Copy()
    $ws.Range($rangeAddr).PasteSpecial($xlPasteFormats) | Out-Null
}

# ====================================================================
# STEP 1: new column K header (2021)
# ====================================================================
Apply-Style "K3" "9"
$ws.Range("K3").Value = 2021

# ====================================================================
# STEP 2: Row 4 (new) = "kg per person" row (previously row 6 content).
# ====================================================================
Apply-Style "A4:C4" "13"
Apply-Style "D4:I4" "22"
Apply-Style "J4" "23"
Apply-Style "K4" "19"

$ws.Range("A4").Value = "Коркунучтуу калдыктардын пайда болушу 1 адамга эсептегенде, килаграмм/адам"
$ws.Range("B4").Value = "Образование опасных отходов в расчете на 1 человека, килограмм/человек"
$ws.Range("C4").Value = "Generation of hazardous waste per person, kilogram / person"
$ws.Range("D4").Value = $row6vals["D"]
$ws.Range("E4").Value = $row6vals["E"]
$ws.Range("F4").Value = $row6vals["F"]
$ws.Range("G4").Value = $row6vals["G"]
$ws.Range("H4").Value = $row6vals["H"]
$ws.Range("I4").Formula = "=I5/I6*1000"
$ws.Range("J4").Value = 1754.6
$ws.Range("K4").Value = 1673.3508218102056
$ws.Rows("4:4").RowHeight = 27

# ====================================================================
# STEP 3: Row 5 (new) = "thousand tons" row (previously row 4 content).
# ====================================================================
Apply-Style "A5:C5" "13"
Apply-Style "D5" "21"
Apply-Style "E5:I5" "22"
Apply-Style "J5:K5" "18"

$ws.Range("A5").Value = "Коркунучтуу калдыктардын пайда болушу (1-3-класстардын коркунучу), миң тонна"
$ws.Range("B5").Value = "Образование опасных отходов (1-3 классов опасности), тыс. тонн"
$ws.Range("C5").Value = "Hazardous waste generation (1-3 hazard classes), thousand tons"
$ws.Range("D5").Value = $row4vals["D"]
$ws.Range("E5").Value = $row4vals["E"]
$ws.Range("F5").Value = $row4vals["F"]
$ws.Range("G5").Value = $row4vals["G"]
$ws.Range("H5").Value = $row4vals["H"]
$ws.Range("I5").Value = $row4vals["I"]
$ws.Range("J5").Value = $row4vals["J"]
$ws.Range("K5").Value = 11290.6
$ws.Rows("5:5").RowHeight = 26.25

# ====================================================================
# STEP 4: Row 6 (new) = "population" row (previously row 5 content).
# ====================================================================
Apply-Style "A6:C6" "14"
Apply-Style "D6:I6" "22"
Apply-Style "J6:K6" "19"

$ws.Range("A6").Value = "Туруктуу калктын саны, миң адам"
$ws.Range("B6").Value = "Численность постоянного населения,  тыс. человек"
$ws.Range("C6").Value = "Resident population, thousand people"
$ws.Range("D6").Value = $row5vals["D"]
$ws.Range("E6").Value = $row5vals["E"]
$ws.Range("F6").Value = $row5vals["F"]
$ws.Range("G6").Value = $row5vals["G"]
$ws.Range("H6").Value = $row5vals["H"]
$ws.Range("I6").Value = $row5vals["I"]
$ws.Range("J6").Value = $row5vals["J"]
$ws.Range("K6").Value = 6747.3
$ws.Rows("6:6").RowHeight = 16.5

# ====================================================================
# STEP 5: Row 7 (new) = "share of neutralized hazardous waste".
# ====================================================================
Apply-Style "A7:C7" "14"
Apply-Style "D7:I7" "22"
Apply-Style "J7:K7" "19"

$ws.Range("A7").Value = "Нейтралдаштырылган кооптуу калдыктардын үлүшү, пайыз"
$ws.Range("B7").Value = "Доля обезвреженных опасных отходов, процентов"
$ws.Range("C7").Value = "Percentage of neutralized hazardous waste, percent"
$ws.Range("D7").Value = ""
$ws.Range("E7").Value = 46.7
$ws.Range("F7").Value = 32.9
$ws.Range("G7").Value = 40.200000000000003
$ws.Range("H7").Value = 38.9
$ws.Range("I7").Value = 53.8
$ws.Range("J7").Value = 56.2
$ws.Range("K7").Value = 57
$ws.Rows("7:7").RowHeight = 15

# ====================================================================
# STEP 6: Row 8 (new) = "share of buried hazardous waste" (bottom,
# bordered row - uses the ORIGINAL bottom-border style family).
# ====================================================================
Apply-Style "A8:C8" "15"
Apply-Style "D8" "16"
Apply-Style "E8:I8" "17"
Apply-Style "J8:K8" "20"

$ws.Range("A8").Value = "Көмүлгөн коркунучтуу калдыктардын үлүшү, пайыз"
$ws.Range("B8").Value = "Доля захороненных опасных отходов, процентов"
$ws.Range("C8").Value = "The share of hazardous waste buried, percent"
$ws.Range("D8").Value = ""
$ws.Range("E8").Value = 0.1
$ws.Range("F8").Value = 0.1
$ws.Range("G8").Value = 0.03
$ws.Range("H8").Value = 0.1
$ws.Range("I8").Value = 0.070000000000000007
$ws.Range("J8").Value = 0
$ws.Range("K8").Value = 0.1
$ws.Rows("8:8").RowHeight = 24.75

# ====================================================================
# STEP 7: clean up scratch cells used for style staging.
# ====================================================================
$ws.Range("N1:N12").Clear()

# ====================================================================
# STEP 8: misc sheet-level tweaks from the diff:
#  - drop the leftover cell selection in the sheet view
#  - bump pageSetup verticalDpi to 300
# ====================================================================
$ws.Range("A1").Select()
$ws.PageSetup.PrintErrors = $ws.PageSetup.PrintErrors
$wb.Save()
